$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend formatting (style) for the new rows 465:491 by copying the format
# of the last existing row (464) down across the new range.
$ws.Range("A464:D464").Copy($ws.Range("A465:D491"))

$data = New-Object "object[,]" 27,4
$data[0,0] = 44539
$data[0,1] = 4
$data[0,2] = 44
$data[0,3] = 534.75935828877
$data[1,0] = 44540
$data[1,1] = 3
$data[1,2] = 34
$data[1,3] = 413.2231404958678
$data[2,0] = 44541
$data[2,1] = 5
$data[2,2] = 38
$data[2,3] = 461.8376276130286
$data[3,0] = 44542
$data[3,1] = 12
$data[3,2] = 48
$data[3,3] = 583.373845405931
$data[4,0] = 44543
$data[4,1] = 11
$data[4,2] = 51
$data[4,3] = 619.8347107438017
$data[5,0] = 44544
$data[5,1] = 5
$data[5,2] = 40
$data[5,3] = 486.1448711716092
$data[6,0] = 44545
$data[6,1] = 0
$data[6,2] = 40
$data[6,3] = 486.1448711716092
$data[7,0] = 44546
$data[7,1] = 8
$data[7,2] = 44
$data[7,3] = 534.75935828877
$data[8,0] = 44547
$data[8,1] = 2
$data[8,2] = 43
$data[8,3] = 522.6057365094798
$data[9,0] = 44548
$data[9,1] = 15
$data[9,2] = 53
$data[9,3] = 644.141954302382
$data[10,0] = 44550
$data[10,1] = 3
$data[10,2] = 44
$data[10,3] = 534.75935828877
$data[11,0] = 44551
$data[11,1] = 6
$data[11,2] = 39
$data[11,3] = 473.9912493923189
$data[12,0] = 44552
$data[12,1] = 0
$data[12,2] = 34
$data[12,3] = 413.2231404958678
$data[13,0] = 44553
$data[13,1] = 4
$data[13,2] = 38
$data[13,3] = 461.8376276130286
$data[14,0] = 44554
$data[14,1] = 0
$data[14,2] = 30
$data[14,3] = 364.6086533787068
$data[15,0] = 44555
$data[15,1] = 4
$data[15,2] = 32
$data[15,3] = 388.9158969372873
$data[16,0] = 44556
$data[16,1] = 3
$data[16,2] = 20
$data[16,3] = 243.0724355858046
$data[17,0] = 44557
$data[17,1] = 2
$data[17,2] = 19
$data[17,3] = 230.9188138065143
$data[18,0] = 44558
$data[18,1] = 7
$data[18,2] = 20
$data[18,3] = 243.0724355858046
$data[19,0] = 44559
$data[19,1] = 1
$data[19,2] = 21
$data[19,3] = 255.2260573650948
$data[20,0] = 44560
$data[20,1] = 2
$data[20,2] = 19
$data[20,3] = 230.9188138065143
$data[21,0] = 44561
$data[21,1] = 7
$data[21,2] = 26
$data[21,3] = 315.994166261546
$data[22,0] = 44562
$data[22,1] = 2
$data[22,2] = 24
$data[22,3] = 291.6869227029655
$data[23,0] = 44563
$data[23,1] = 11
$data[23,2] = 32
$data[23,3] = 388.9158969372873
$data[24,0] = 44564
$data[24,1] = 10
$data[24,2] = 40
$data[24,3] = 486.1448711716092
$data[25,0] = 44565
$data[25,1] = 6
$data[25,2] = 39
$data[25,3] = 473.9912493923189
$data[26,0] = 44566
$data[26,1] = 20
$data[26,2] = 58
$data[26,3] = 704.9100631988332

$ws.Range("A465:D491").Value = $data
